$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("control_panel")

# Clear the values (but keep formatting/style) for cells C10:I10
$ws.Range("C10:I10").ClearContents()

# Move the active selection from I10 to C10
$ws.Activate()
$ws.Range("C10").Select()
